$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23, shifting existing rows 23:129 down to 24:130.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23. Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the
# same values the (old) row 23 had; D,J,K,L,M,P get new values per the diff.
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "Macroferia Regional de Talca"
$ws.Range("C23").Value = "Maule"
$ws.Range("D23").Value = "2021-12-28"
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = 100112031
$ws.Range("G23").Value = "Poroto verde"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 33000
$ws.Range("L23").Value = 33000
$ws.Range("M23").Value = 33000
$ws.Range("N23").Value = "$/saco 25 kilos"
$ws.Range("O23").Value = "Región del Maule"
$ws.Range("P23").Value = 1320
$ws.Range("Q23").Value = 25
$ws.Range("R23").Value = "Hortaliza"
